# Handles float input without breaking stuff
#
# The generated "Marksheet" worksheet previously treated the student as
# fully absent (every question blank / "Not Attempt", summary = 0/0,
# grand total "Absent"). This updates it to reflect the student's real
# (floating point capable) submitted answers: the summary scores are
# recomputed, each attempted question in the first answer block (column
# A) is filled in and marked correct/incorrect via the matching named
# style, and the now-unused third answer block (columns G:H) plus the
# unused tail of the second answer block (D19:E40) are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-Style($fromRef, $toRef) {
    $ws.Range($fromRef).Copy() | Out-Null
    $ws.Range($toRef).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

function Set-Answer($cellRef, $styleSourceRef, $text) {
    Copy-Style $styleSourceRef $cellRef
    $ws.Range($cellRef).Value = $text
}

# ---------------------------------------------------------------
# 1. Summary block (rows 10-12): header styling + recomputed scores
# ---------------------------------------------------------------
Copy-Style "A9" "A10"
Copy-Style "A9" "A11"
Copy-Style "A9" "A12"

$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 7
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 36
$ws.Range("C12").Value = -7
$ws.Range("E12").Value = "29/112"

# ---------------------------------------------------------------
# 2. Remove the third answer block (columns G:H, rows 15-21) and
#    the unused tail of the second answer block (D19:E40)
# ---------------------------------------------------------------
$ws.Range("G15:H21").Clear()
$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------
# 3. Fill in the student's actual answers for the first answer
#    block (column A): correctStyle / incorrectStyle on top of the
#    existing normalStyle ("not attempted") blanks
# ---------------------------------------------------------------
Set-Answer "A18" "B10" "Option B"
Set-Answer "A19" "B10" "Option C"
Set-Answer "A21" "B10" "Option C"
Set-Answer "A22" "C10" "Option A"
Set-Answer "A23" "B10" "Option D"
Set-Answer "A24" "C10" "Option C"
Set-Answer "A26" "C10" "Option D"
Set-Answer "A28" "C10" "Option B"
Set-Answer "A29" "B10" "Option D"
Set-Answer "A33" "B10" "Option D"
Set-Answer "A35" "B10" "Option D"
Set-Answer "A37" "C10" "Option B"
Set-Answer "A39" "B10" "Option D"
Set-Answer "A40" "C10" "Option B"

# ---------------------------------------------------------------
# 4. Two cells in the (remaining) second answer block also change
# ---------------------------------------------------------------
Set-Answer "D16" "C10" "Option C"
Set-Answer "D18" "B10" "Option D"

Write-Host "Done. UsedRange:" ($ws.UsedRange.Address())
